$d = $word.ActiveDocument
# Remove the old _GoBack bookmark first
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Now insert the new run + bookmark into the first paragraph
$r = $d.Range(0, 12)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="000B0A8E"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Multibinding</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>, IMultiValueConverter</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Output $d.Paragraphs(1).Range.Text
Write-Output $d.Paragraphs.Count
Write-Output $d.Bookmarks.Exists("_GoBack")
